$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")

# Update header row of Hoja1: "Tratamiento"/"Bloque"/"H"/"RS" -> "Trat"/"rep"/"H"/"rend seco"
$ws1.Range("A1").Value = "Trat"
$ws1.Range("B1").Value = "rep"
$ws1.Range("C1").Value = "H"
$ws1.Range("D1").Value = "rend seco"

# Make Hoja1 the active sheet/tab (was "prom")
$ws1.Activate()

# Update the selection on Hoja1 to D1 (was C1:D1)
$ws1.Range("D1").Select()
